# Update "Hortaliza, Femacal de La Calera - Pepino dulce" data.
# The weekly records for rows 5-6 are swapped with rows 9-10, and the
# weekly records for rows 7-8 are swapped with rows 11-12 (columns
# D, J, K, L, M, N, P, Q only - the other columns are identical across
# the affected rows already).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters used per row: D=4, J=10, K=11, L=12, M=13, N=14, P=16, Q=17
$cols = @(4, 10, 11, 12, 13, 14, 16, 17)

# Target values, keyed by row number, after the change.
$values = @{
    5  = @{ 4 = 44242; 10 = 60; 11 = 13000; 12 = 13000; 13 = 13000; 14 = "`$/bandeja 18 kilos"; 16 = 722;  17 = 18 }
    6  = @{ 4 = 44242; 10 = 50; 11 = 10000; 12 = 10000; 13 = 10000; 14 = "`$/bandeja 18 kilos"; 16 = 556;  17 = 18 }
    7  = @{ 4 = 44424; 10 = 75; 11 = 18000; 12 = 18000; 13 = 18000; 14 = "`$/caja 15 kilos";    16 = 1200; 17 = 15 }
    8  = @{ 4 = 44424; 10 = 50; 11 = 12000; 12 = 12000; 13 = 12000; 14 = "`$/caja 15 kilos";    16 = 800;  17 = 15 }
    9  = @{ 4 = 44238; 10 = 90; 11 = 13000; 12 = 13000; 13 = 13000; 14 = "`$/bandeja 18 kilos"; 16 = 722;  17 = 18 }
    10 = @{ 4 = 44238; 10 = 80; 11 = 11000; 12 = 11000; 13 = 11000; 14 = "`$/bandeja 18 kilos"; 16 = 611;  17 = 18 }
    11 = @{ 4 = 44536; 10 = 87; 11 = 22000; 12 = 22000; 13 = 22000; 14 = "`$/bandeja 18 kilos"; 16 = 1222; 17 = 18 }
    12 = @{ 4 = 44536; 10 = 80; 11 = 20000; 12 = 20000; 13 = 20000; 14 = "`$/bandeja 18 kilos"; 16 = 1111; 17 = 18 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    foreach ($col in $cols) {
        $ws.Cells.Item($row, $col).Value = $rowVals[$col]
    }
}
